$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.760.09"
$ws.Range("E2").Value = "  -1.69%  "
$ws.Range("D3").Value = "'1.895.52"
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("D5").Value = "'311.93"
$ws.Range("E5").Value = "  -1.49%  "
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("D7").Value = "'0.4926"
$ws.Range("E7").Value = "  +1.63%  "
$ws.Range("D8").Value = "'0.3791"
$ws.Range("E8").Value = "  -1.57%  "
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("D10").Value = "'0.9107"
$ws.Range("E10").Value = "  -4.17%  "
$ws.Range("D11").Value = "'20.63"
$ws.Range("E11").Value = "  -1.51%  "
$ws.Range("D12").Value = "'0.07623"
$ws.Range("E12").Value = "  -2.21%  "
$ws.Range("D13").Value = "'1.900.71"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("D14").Value = "'5.460"
$ws.Range("E14").Value = "  -1.74%  "
$ws.Range("D15").Value = "'6.647"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").Value = "'91.10"
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").Value = "'0.000008731"
$ws.Range("E18").Value = "  -1.53%  "
$ws.Range("D19").Value = "'0.9995"
$ws.Range("D20").Value = "'27.742.47"
$ws.Range("E20").Value = "  -1.77%  "
$ws.Range("D21").Value = "'14.47"
$ws.Range("E21").Value = "  -3.73%  "
$ws.Range("D22").Value = "'5.115"
$ws.Range("E22").Value = "  -1.15%  "
$ws.Range("D23").Value = "'2.114.03"
$ws.Range("E23").Value = "  -2.55%  "
$ws.Range("D24").Value = "'10.75"
$ws.Range("E24").Value = "  -1.89%  "
$ws.Range("D25").Value = "'153.92"
$ws.Range("E25").Value = "  -1.46%  "
$ws.Range("D26").Value = "'1.848"
$ws.Range("E26").Value = "  -4.26%  "
$ws.Range("D27").Value = "'2.177"
$ws.Range("E27").Value = "  +3.01%  "
$ws.Range("E28").Value = "  -1.45%  "
$ws.Range("D29").Value = "'115.25"
$ws.Range("E29").Value = "  -1.82%  "
$ws.Range("D30").Value = "'4.880"
$ws.Range("E30").Value = "  -2.55%  "
$ws.Range("D31").Value = "'0.08934"
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("E32").Value = "  -4.03%  "
$ws.Range("D33").Value = "'1.229"
$ws.Range("E33").Value = "  -1.81%  "
$ws.Range("D34").Value = "'0.7667"
$ws.Range("D35").Value = "'4.648"
$ws.Range("E35").Value = "  -0.86%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.02045"
$ws.Range("E36").Value = "  -0.59%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "'2.558"
$ws.Range("E37").Value = "  -8.05%  "
$ws.Range("D38").Value = "'1.099"
$ws.Range("E38").Value = "  -2.86%  "
$ws.Range("D39").Value = "'0.5495"
$ws.Range("E39").Value = "  -1.83%  "
$ws.Range("D40").Value = "'0.05282"
$ws.Range("E40").Value = "  -1.83%  "
$ws.Range("D41").Value = "'2.989"
$ws.Range("E41").Value = "  -1.60%  "
$ws.Range("D42").Value = "'6.898"
$ws.Range("E42").Value = "  -3.18%  "
$ws.Range("D43").Value = "'8.563"
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("D44").Value = "'112.63"
$ws.Range("E44").Value = "  +4.77%  "
$ws.Range("D45").Value = "'0.1520"
$ws.Range("E45").Value = "  -1.10%  "
$ws.Range("D46").Value = "'10.56"
$ws.Range("E46").Value = "  -2.78%  "
$ws.Range("D47").Value = "'0.4793"
$ws.Range("E47").Value = "  -2.39%  "
$ws.Range("D48").Value = "'0.9990"
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("D49").Value = "'1.632"
$ws.Range("E49").Value = "  -2.77%  "
$ws.Range("D50").Value = "'67.44"
$ws.Range("E50").Value = "  -2.95%  "
$ws.Range("D51").Value = "'0.06059"
$ws.Range("E51").Value = "  -1.53%  "
